# Sync attendance_reports - fix "Recorded By" (column G) ordering.
# The first two comma-separated entries in several "Recorded By" values
# were recorded in the wrong order (e.g. "dnasr281@gmail.com, System"
# should read "System, dnasr281@gmail.com"). This swaps the first two
# tokens for the three known mis-ordered patterns, leaving already
# correctly ordered / single-entry values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact (old -> new) replacements observed in the sheet's "Recorded By" column.
$map = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
